$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '55.968.90'
$ws.Range("E2").Value = '  +3.17%  '

$ws.Range("D3").Value = '2.310.28'
$ws.Range("E3").Value = '  +1.99%  '

$ws.Range("E4").Value = '  +0.22%  '

$ws.Range("D5").Value = '516.04'
$ws.Range("E5").Value = '  +3.74%  '

$ws.Range("D6").Value = '132.63'
$ws.Range("E6").Value = '  +2.83%  '

$ws.Range("D7").Value = '0.995'
$ws.Range("E7").Value = '  -0.33%  '

$ws.Range("D8").Value = '0.533'
$ws.Range("E8").Value = '  +1.57%  '

$ws.Range("D9").Value = '2.325.84'
$ws.Range("E9").Value = '  +2.39%  '

$ws.Range("D10").Value = '0.104'
$ws.Range("E10").Value = '  +8.64%  '

$ws.Range("E11").Value = '  +0.74%  '

$ws.Range("D12").Value = '5.16'
$ws.Range("E12").Value = '  +8.09%  '

$ws.Range("E13").Value = '  +1.59%  '

$ws.Range("D14").Value = '24.02'
$ws.Range("E14").Value = '  +4.90%  '

$ws.Range("D15").Value = '2.723.47'
$ws.Range("E15").Value = '  +2.13%  '

$ws.Range("D16").Value = '56.195.07'
$ws.Range("E16").Value = '  +3.64%  '

$ws.Range("E17").Value = '  +4.36%  '

$ws.Range("D18").Value = '2.301.79'
$ws.Range("E18").Value = '  +1.91%  '

$ws.Range("D19").Value = '10.51'
$ws.Range("E19").Value = '  +2.69%  '

$ws.Range("D20").Value = '4.27'
$ws.Range("E20").Value = '  +2.96%  '

$ws.Range("D21").Value = '321.72'
$ws.Range("E21").Value = '  +5.97%  '

$ws.Range("D22").Value = '6.66'
$ws.Range("E22").Value = '  +5.57%  '

$ws.Range("D23").Value = '0.999'
$ws.Range("E23").Value = '  -0.03%  '

$ws.Range("D24").Value = '60.66'
$ws.Range("E24").Value = '  -0.48%  '

$ws.Range("D25").Value = '0.990'
$ws.Range("E25").Value = '  -0.46%  '

$ws.Range("E26").Value = '  +5.96%  '

$ws.Range("E27").Value = '  +4.16%  '

$ws.Range("D28").Value = '171.72'
$ws.Range("E28").Value = '  +0.16%  '

$ws.Range("E29").Value = '  +10.31%  '

$ws.Range("D30").Value = '0.0₃0726'
$ws.Range("E30").Value = '  +5.00%  '

$ws.Range("B31").Value = 'Aptos'
$ws.Range("C31").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D31").Value = '6.25'
$ws.Range("E31").Value = '  +5.72%  '

$ws.Range("B32").Value = 'PancakeSwap'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D32").Value = '1.68'
$ws.Range("E32").Value = '  +4.24%  '

$ws.Range("D33").Value = '18.33'
$ws.Range("E33").Value = '  +3.20%  '

$ws.Range("E34").Value = '  +0.01%  '

$ws.Range("E35").Value = '  -0.41%  '

$ws.Range("E36").Value = '  +6.04%  '

$ws.Range("D37").Value = '0.928'
$ws.Range("E37").Value = '  -0.83%  '

$ws.Range("E38").Value = '  +7.62%  '

$ws.Range("E39").Value = '  +8.37%  '

$ws.Range("D40").Value = '37.28'
$ws.Range("E40").Value = '  +3.72%  '

$ws.Range("D41").Value = '0.382'
$ws.Range("E41").Value = '  +2.01%  '

$ws.Range("B42").Value = 'Aave'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D42").Value = '139.64'
$ws.Range("E42").Value = '  +11.73%  '

$ws.Range("B43").Value = 'Filecoin'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D43").Value = '3.57'
$ws.Range("E43").Value = '  +6.14%  '

$ws.Range("E44").Value = '  +6.60%  '

$ws.Range("D45").Value = '269.88'
$ws.Range("E45").Value = '  +11.85%  '

$ws.Range("E46").Value = '  +3.40%  '

$ws.Range("D47").Value = '0.0926'

$ws.Range("D48").Value = '0.552'
$ws.Range("E48").Value = '  +1.11%  '

$ws.Range("E49").Value = '  +2.26%  '

$ws.Range("E50").Value = '  +4.92%  '

$ws.Range("D51").Value = '16.91'
$ws.Range("E51").Value = '  +4.84%  '

